$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(2833958, 27.7, 785395, 4,   113896, 3.9, 111194, 1.4, 40407, 62.9, 1783065, 37.1, 1050893)
    3  = @(2471589, 24.4, 603451, 4,   97997,  4,   99845,  1,   25528, 66.5, 1644768, 33.5, 826821)
    4  = @(226810,  51.4, 116675, 3.3, 7382,   3.4, 7675,   5.2, 11695, 36.8, 83383,   63.2, 143427)
    5  = @(113405,  45.8, 51983,  7.1, 8099,   3.1, 3469,   1.8, 2034,  42.2, 47820,   57.8, 65585)
    6  = @(22154,   60,   13287,  1.9, 418,    0.9, 205,    5.2, 1150,  32,   7094,    68,   15060)
    7  = @(1251549, 26.8, 335825, 4.4, 54947,  3.9, 48329,  1.3, 15870, 63.6, 796577,  36.4, 454972)
    8  = @(1582409, 28.6, 451938, 3.6, 56971,  4.1, 65185,  1.6, 24672, 62.2, 983644,  37.8, 598765)
    9  = @(290232,  26.8, 77773,  0.8, 2364,   2.4, 6904,   1.5, 4234,  68.59999999999999, 198957, 31.4, 91274)
    10 = @(1089984, 22.2, 242310, 0.9, 10040,  5.1, 55205,  1.8, 19733, 70,   762696,  30,   327288)
    11 = @(871987,  26.4, 229861, 3,   26406,  3.3, 29001,  1.4, 12007, 65.90000000000001, 574713, 34.1, 297274)
    12 = @(653990,  40.5, 264609, 12,  78244,  3.1, 20489,  0.7, 4536,  43.7, 286113,  56.3, 367878)
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O")

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $ws.Range("$col$rowNum").Value = $values[$i]
    }
}
